$d = $word.ActiveDocument

# 1. Remove the word "Solo " (with trailing space) that precedes "Projects"
#    in the "Solo Projects" heading - it doesn't read well to ATS parsers.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Solo Projects", $true, $false, $false, $false, $false, $true, 1, $false, "Projects", 2)

# 2. Word leaves its "last edit" (_GoBack) bookmark behind wherever the
#    cursor ends up after the edit session. Re-create that bookmark at
#    its resulting location, inside "component" (between "comp" and
#    "onent"), splitting the run there - matching the authored change.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found = $find2.Execute("component", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPos = $find2.Parent.Start + 4
    $bmRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
